# Apply the changes described by the commit diff:
#   1. Fix the typo "alaxa" -> "alexa" in cell A1 of Sheet1.
#   2. Give cell B2 a thin black bottom border (new border record + cellXf
#      pointing at it with applyBorder="1").
#   3. Leave the selection on A1 (the sheet's default cell) instead of B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the misspelled shared string in A1 ---------------------
$ws.Range("A1").Value = "alexa"

# --- 2. Add a thin bottom border to B2 ----------------------------------
$bottomBorder = $ws.Range("B2").Borders.Item(9)   # 9 = xlEdgeBottom
$bottomBorder.Color = 0                            # black
$bottomBorder.LineStyle = 1                        # 1 = xlContinuous (thin)

# --- 3. Return the active selection to A1 -------------------------------
$ws.Range("A1").Select()
